# rec_selectivity_coastwide_cdf_star: refresh C_l (col B) and q (col D)
# for Sheet1 rows 2-27 (l_in_bin 6..31). Columns A (l_in_bin) and C (E) are
# untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# sheet row number -> @(new C_l, new q), values as strings so the
# "e-NN" scientific-notation literals parse; cast to [double] on write.
$newValues = [ordered]@{
    "2" = @("1004.5607611341402", "7.4894104806766393e-11")
    "3" = @("16345.936218261719", "9.9284851173830901e-11")
    "4" = @("115572.06396484375", "5.8319815643415041e-10")
    "5" = @("573382.234375", "2.5549280469050473e-09")
    "6" = @("1461357.8125", "1.6636956701177041e-08")
    "7" = @("2318790.75", "2.7467175200968086e-08")
    "8" = @("3563264.3125", "5.7023243016374181e-08")
    "9" = @("4596000", "5.467463282116114e-08")
    "10" = @("3742583", "3.3140537425424554e-08")
    "11" = @("4707804", "2.4003771414982111e-08")
    "12" = @("2927325.375", "3.0520627802843592e-08")
    "13" = @("2233480.3125", "2.5553736904271318e-08")
    "14" = @("1100141.5", "2.2067942140324703e-08")
    "15" = @("1141520.875", "2.3915026403642514e-08")
    "16" = @("488117.171875", "1.2478896138645723e-08")
    "17" = @("442970.9375", "2.1153725882072649e-08")
    "18" = @("164226.234375", "5.606165753846426e-09")
    "19" = @("98302.7109375", "4.7247077361589618e-09")
    "20" = @("56294.08984375", "2.6261053331921858e-09")
    "21" = @("30841.416015625", "3.1530840161053675e-09")
    "22" = @("16165.15380859375", "1.6372625477600877e-09")
    "23" = @("8105.861328125", "1.7877554991940769e-09")
    "24" = @("2308.74560546875", "5.0290921427276203e-10")
    "25" = @("2648.979736328125", "1.6847877537529143e-09")
    "26" = @("1189.1361389160156", "1.9266994666367054e-09")
    "27" = @("310.00677490234375", "7.5343342675893155e-10")
}

foreach ($row in $newValues.Keys) {
    $pair = $newValues[$row]
    $ws.Cells.Item([int]$row, 2).Value = [double]$pair[0]   # column B: C_l
    $ws.Cells.Item([int]$row, 4).Value = [double]$pair[1]   # column D: q
}
